$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.02678108215332
$ws.Range("B1").Value = 3.31135892868042
$ws.Range("C1").Value = 3.58078145980835
$ws.Range("D1").Value = 2.049391269683838
$ws.Range("E1").Value = 1.177251100540161
